$d = $word.ActiveDocument

# "Salva" + " tra uno step " -> "Condividere dati tra uno step "
$d.Content.Find.Execute("Salva tra uno step ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Condividere dati tra uno step ", 2) | Out-Null

# " successivo il numero generato." -> " successivo."
$d.Content.Find.Execute(" successivo il numero generato.", $true, $false, $false, $false, $false, `
    $true, 1, $false, " successivo.", 2) | Out-Null

# "Generatore di numeri casuali" -> "Generatore di numeri sequenziali"
$d.Content.Find.Execute("Generatore di numeri casuali", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Generatore di numeri sequenziali", 2) | Out-Null

# Move the "_GoBack" bookmark so that it now sits right after "sequenziali"
# (end of that paragraph) instead of after "Impl" further down the document.
# A collapsed bookmark range landing exactly at (paragraph end - 1) is mishandled
# by the engine, so we work around it: temporarily insert a placeholder run,
# anchor the bookmark next to it, then delete the placeholder again.
$gp = $d.Content.Find
$gp.Execute("Generatore di numeri sequenziali", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$genPara = $gp.Parent.Paragraphs(1)
$paraEnd = $genPara.Range.End

$placeholder = $d.Range($paraEnd - 1, $paraEnd - 1)
$placeholder.InsertBefore("@@MARK@@")

$newEnd = $genPara.Range.End
$bmRange = $d.Range($newEnd - 9, $newEnd - 9)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholderRange = $d.Range($newEnd - 9, $newEnd - 1)
$placeholderRange.Text = ""

# "IMPL:" -> "COMMON:"
$d.Content.Find.Execute("IMPL:", $true, $false, $false, $false, $false, `
    $true, 1, $false, "COMMON:", 2) | Out-Null
